# Add a new "Tags" column to the accessioning template, right after the
# "Specimen Notes" column (AM) and before the "Files" column (old AN, new AO).
#
# This mirrors a user selecting column AN and doing Insert (Shift cells right),
# which pushes all the metadata columns (Files, Case Files, Genome Build,
# Variant Type, BAM Sample ID, ...) one column over, then typing the new
# header "Tags" in row 2 and the example text "Comma-separated" in row 3
# (consistent with the neighboring "Files" column which is also
# comma-separated), and finally landing the selection down at AN7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at AN (column 40); this shifts the previous
# AN:CC columns (and their column-width/style definitions) one slot to the
# right, exactly like Excel's "Insert Sheet Columns" command.
$ws.Columns("AN:AN").Insert()

# The newly inserted column otherwise keeps Excel's generic default width;
# give it the same width as its left neighbour (AM), matching the look of
# the rest of the "Specimen Information" block.
$ws.Columns("AN:AN").ColumnWidth = $ws.Columns("AM:AM").ColumnWidth

# New column header (row 2) and example/description text (row 3) for the
# newly created "Tags" column.
$ws.Range("AN2").Value = "Tags"
$ws.Range("AN3").Value = "Comma-separated"

# Leave the selection where the editor ended up after making the change.
$ws.Range("AN7").Select()
